# Natmi following Dr Hou advice
# Re-run of the Cd38-Pecam1 ligand/receptor analysis for YoungD0 with an
# added "sCs" cluster. This rewrites the existing result rows (2-7) with
# their recomputed statistics and appends the new rows (8-10) produced by
# including "sCs" as a sending cluster.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ A="ECs"; B="Cd38"; C="Pecam1"; D="ECs"; E=3; F=1; G=40.61064433333333; H=121.831933; I=0.9365419382944963; J=0.9365419382944962; K=3; L=1; M=234.5813823333333; N=703.744147; O=0.9782746109134588; P=0.9782746109134588; Q=9526.501085160682; R=85738.50976644615; S=0.9161952002891849; T=0.9161952002891848 },
    @{ A="ECs"; B="Cd38"; C="Pecam1"; D="FAPs"; E=3; F=1; G=40.61064433333333; H=121.831933; I=0.9365419382944963; J=0.9365419382944962; K=3; L=1; M=2.110333333333333; N=6.331; O=0.008800721950008783; P=0.008800721950008783; Q=85.70199642477778; R=771.317967823; S=0.008242245193452144; T=0.008242245193452144 },
    @{ A="ECs"; B="Cd38"; C="Pecam1"; D="sCs"; E=3; F=1; G=40.61064433333333; H=121.831933; I=0.9365419382944963; J=0.9365419382944962; K=3; L=1; M=3.099218; N=9.297654; O=0.01292466713653245; P=0.01292466713653245; Q=125.8612399094647; R=1132.751159185182; S=0.01210449281185928; T=0.01210449281185928 },
    @{ A="FAPs"; B="Cd38"; C="Pecam1"; D="ECs"; E=3; F=1; G=2.454549333333333; H=7.363648; I=0.05660556309845621; J=0.05660556309845621; K=3; L=1; M=234.5813823333333; N=703.744147; O=0.9782746109134588; P=0.9782746109134588; Q=575.791575618695; R=5182.124180568256; S=0.05537578521567949; T=0.05537578521567948 },
    @{ A="FAPs"; B="Cd38"; C="Pecam1"; D="FAPs"; E=3; F=1; G=2.454549333333333; H=7.363648; I=0.05660556309845621; J=0.05660556309845621; K=3; L=1; M=2.110333333333333; N=6.331; O=0.008800721950008783; P=0.008800721950008783; Q=5.179917276444443; R=46.61925548799999; S=0.0004981698216531907; T=0.0004981698216531907 },
    @{ A="FAPs"; B="Cd38"; C="Pecam1"; D="sCs"; E=3; F=1; G=2.454549333333333; H=7.363648; I=0.05660556309845621; J=0.05660556309845621; K=3; L=1; M=3.099218; N=9.297654; O=0.01292466713653245; P=0.01292466713653245; Q=7.607183475754666; R=68.46465128179199; S=0.0007316080611235312; T=0.0007316080611235311 },
    @{ A="sCs"; B="Cd38"; C="Pecam1"; D="ECs"; E=2; F=0.6666666666666666; G=0.2971403333333333; H=0.891421; I=0.006852498607047613; J=0.006852498607047612; K=3; L=1; M=234.5813823333333; N=703.744147; O=0.9782746109134588; P=0.9782746109134588; Q=69.70359014032077; R=627.3323112628871; S=0.006703625408594522; T=0.006703625408594521 },
    @{ A="sCs"; B="Cd38"; C="Pecam1"; D="FAPs"; E=2; F=0.6666666666666666; G=0.2971403333333333; H=0.891421; I=0.006852498607047613; J=0.006852498607047612; K=3; L=1; M=2.110333333333333; N=6.331; O=0.008800721950008783; P=0.008800721950008783; Q=0.6270651501111111; R=5.643586351; S=0.00006030693490344854; T=0.00006030693490344853 },
    @{ A="sCs"; B="Cd38"; C="Pecam1"; D="sCs"; E=2; F=0.6666666666666666; G=0.2971403333333333; H=0.891421; I=0.006852498607047613; J=0.006852498607047612; K=3; L=1; M=3.099218; N=9.297654; O=0.01292466713653245; P=0.01292466713653245; Q=0.9209026695926666; R=8.288124026334; S=0.0000885662635496427; T=0.00008856626354964269 }
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$r = 2
foreach ($row in $rows) {
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value = $row[$col]
    }
    $r++
}